# Apply 2024-10-18 violent crime data update across all affected worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 11).Value = 5882
$ws.Cells.Item(3, 11).Value = 6057
$ws.Cells.Item(4, 7).Value = 1489
$ws.Cells.Item(4, 11).Value = 1260
$ws.Cells.Item(5, 11).Value = 431
$ws.Cells.Item(6, 11).Value = 6663
$ws.Cells.Item(7, 7).Value = 24717
$ws.Cells.Item(7, 11).Value = 20293

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 11).Value = 142
$ws.Cells.Item(7, 11).Value = 589
$ws.Cells.Item(8, 11).Value = 1343
$ws.Cells.Item(9, 11).Value = 90
$ws.Cells.Item(11, 11).Value = 385
$ws.Cells.Item(12, 11).Value = 36
$ws.Cells.Item(14, 11).Value = 104
$ws.Cells.Item(19, 11).Value = 584
$ws.Cells.Item(20, 11).Value = 478
$ws.Cells.Item(27, 11).Value = 189
$ws.Cells.Item(29, 11).Value = 1107
$ws.Cells.Item(30, 11).Value = 77
$ws.Cells.Item(31, 11).Value = 226
$ws.Cells.Item(33, 11).Value = 875
$ws.Cells.Item(36, 11).Value = 258
$ws.Cells.Item(37, 11).Value = 687
$ws.Cells.Item(41, 11).Value = 141
$ws.Cells.Item(42, 11).Value = 754
$ws.Cells.Item(45, 11).Value = 26
$ws.Cells.Item(47, 11).Value = 141
$ws.Cells.Item(49, 11).Value = 110
$ws.Cells.Item(51, 11).Value = 259
$ws.Cells.Item(52, 11).Value = 536
$ws.Cells.Item(53, 11).Value = 259
$ws.Cells.Item(54, 11).Value = 394
$ws.Cells.Item(55, 11).Value = 223
$ws.Cells.Item(63, 7).Value = 291
$ws.Cells.Item(65, 11).Value = 473
$ws.Cells.Item(67, 11).Value = 795
$ws.Cells.Item(68, 11).Value = 54
$ws.Cells.Item(72, 11).Value = 97
$ws.Cells.Item(73, 11).Value = 179
$ws.Cells.Item(76, 11).Value = 276
$ws.Cells.Item(77, 11).Value = 143
$ws.Cells.Item(78, 11).Value = 228
$ws.Cells.Item(79, 11).Value = 502
$ws.Cells.Item(80, 11).Value = 71
$ws.Cells.Item(81, 11).Value = 15
$ws.Cells.Item(83, 11).Value = 453
$ws.Cells.Item(84, 11).Value = 160
$ws.Cells.Item(85, 11).Value = 949
$ws.Cells.Item(88, 11).Value = 219
$ws.Cells.Item(90, 11).Value = 187
$ws.Cells.Item(91, 11).Value = 230
$ws.Cells.Item(94, 11).Value = 272
$ws.Cells.Item(97, 11).Value = 161
$ws.Cells.Item(98, 11).Value = 97
$ws.Cells.Item(101, 7).Value = 24717
$ws.Cells.Item(101, 11).Value = 20293

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 104

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 11).Value = 197
$ws.Cells.Item(3, 11).Value = 192
$ws.Cells.Item(6, 11).Value = 157
$ws.Cells.Item(7, 11).Value = 589

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 11).Value = 137
$ws.Cells.Item(3, 11).Value = 99
$ws.Cells.Item(7, 11).Value = 385

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(2, 11).Value = 309
$ws.Cells.Item(3, 11).Value = 325
$ws.Cells.Item(6, 11).Value = 235
$ws.Cells.Item(7, 11).Value = 949

$ws = $wb.Worksheets.Item(9)
$ws.Cells.Item(2, 11).Value = 144
$ws.Cells.Item(3, 11).Value = 156
$ws.Cells.Item(5, 11).Value = 19
$ws.Cells.Item(7, 11).Value = 536

$ws = $wb.Worksheets.Item(11)
$ws.Cells.Item(2, 11).Value = 66
$ws.Cells.Item(3, 11).Value = 68
$ws.Cells.Item(6, 11).Value = 112
$ws.Cells.Item(7, 11).Value = 259

$ws = $wb.Worksheets.Item(12)
$ws.Cells.Item(2, 11).Value = 372
$ws.Cells.Item(3, 11).Value = 411
$ws.Cells.Item(5, 11).Value = 39
$ws.Cells.Item(7, 11).Value = 1343

$ws = $wb.Worksheets.Item(13)
$ws.Cells.Item(3, 11).Value = 161
$ws.Cells.Item(6, 11).Value = 104
$ws.Cells.Item(7, 11).Value = 453

$ws = $wb.Worksheets.Item(14)
$ws.Cells.Item(4, 11).Value = 42
$ws.Cells.Item(5, 11).Value = 22
$ws.Cells.Item(7, 11).Value = 875

$ws = $wb.Worksheets.Item(16)
$ws.Cells.Item(3, 11).Value = 226
$ws.Cells.Item(4, 11).Value = 33
$ws.Cells.Item(6, 11).Value = 200
$ws.Cells.Item(7, 11).Value = 687

$ws = $wb.Worksheets.Item(17)
$ws.Cells.Item(2, 11).Value = 153
$ws.Cells.Item(3, 11).Value = 116
$ws.Cells.Item(6, 11).Value = 175
$ws.Cells.Item(7, 11).Value = 473

$ws = $wb.Worksheets.Item(19)
$ws.Cells.Item(2, 11).Value = 21
$ws.Cells.Item(7, 11).Value = 77

$ws = $wb.Worksheets.Item(20)
$ws.Cells.Item(3, 11).Value = 61
$ws.Cells.Item(7, 11).Value = 226

$ws = $wb.Worksheets.Item(21)
$ws.Cells.Item(2, 11).Value = 222
$ws.Cells.Item(3, 11).Value = 283
$ws.Cells.Item(6, 11).Value = 227
$ws.Cells.Item(7, 11).Value = 795

$ws = $wb.Worksheets.Item(22)
$ws.Cells.Item(3, 11).Value = 61
$ws.Cells.Item(7, 11).Value = 160

$ws = $wb.Worksheets.Item(23)
$ws.Cells.Item(3, 11).Value = 22
$ws.Cells.Item(7, 11).Value = 110

$ws = $wb.Worksheets.Item(24)
$ws.Cells.Item(6, 11).Value = 212
$ws.Cells.Item(7, 11).Value = 394

$ws = $wb.Worksheets.Item(25)
$ws.Cells.Item(3, 11).Value = 399
$ws.Cells.Item(6, 11).Value = 313
$ws.Cells.Item(7, 11).Value = 1107

$ws = $wb.Worksheets.Item(27)
$ws.Cells.Item(3, 11).Value = 179
$ws.Cells.Item(7, 11).Value = 584

$ws = $wb.Worksheets.Item(29)
$ws.Cells.Item(6, 11).Value = 144
$ws.Cells.Item(7, 11).Value = 276

$ws = $wb.Worksheets.Item(30)
$ws.Cells.Item(4, 11).Value = 7
$ws.Cells.Item(7, 11).Value = 142

$ws = $wb.Worksheets.Item(31)
$ws.Cells.Item(4, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 141

$ws = $wb.Worksheets.Item(32)
$ws.Cells.Item(2, 11).Value = 205
$ws.Cells.Item(3, 11).Value = 231
$ws.Cells.Item(4, 11).Value = 30
$ws.Cells.Item(6, 11).Value = 280
$ws.Cells.Item(7, 11).Value = 754

$ws = $wb.Worksheets.Item(35)
$ws.Cells.Item(6, 11).Value = 81
$ws.Cells.Item(7, 11).Value = 228

$ws = $wb.Worksheets.Item(36)
$ws.Cells.Item(4, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 223

$ws = $wb.Worksheets.Item(40)
$ws.Cells.Item(5, 11).Value = 5
$ws.Cells.Item(7, 11).Value = 230

$ws = $wb.Worksheets.Item(42)
$ws.Cells.Item(2, 11).Value = 170
$ws.Cells.Item(3, 11).Value = 162
$ws.Cells.Item(7, 11).Value = 502

$ws = $wb.Worksheets.Item(44)
$ws.Cells.Item(2, 11).Value = 159
$ws.Cells.Item(6, 11).Value = 136
$ws.Cells.Item(7, 11).Value = 478

$ws = $wb.Worksheets.Item(47)
$ws.Cells.Item(2, 11).Value = 102
$ws.Cells.Item(7, 11).Value = 258

$ws = $wb.Worksheets.Item(51)
$ws.Cells.Item(2, 11).Value = 72
$ws.Cells.Item(6, 11).Value = 121
$ws.Cells.Item(7, 11).Value = 272

$ws = $wb.Worksheets.Item(53)
$ws.Cells.Item(2, 11).Value = 41
$ws.Cells.Item(7, 11).Value = 141

$ws = $wb.Worksheets.Item(55)
$ws.Cells.Item(2, 11).Value = 18
$ws.Cells.Item(6, 11).Value = 55
$ws.Cells.Item(7, 11).Value = 97

$ws = $wb.Worksheets.Item(61)
$ws.Cells.Item(3, 11).Value = 33
$ws.Cells.Item(7, 11).Value = 90

$ws = $wb.Worksheets.Item(62)
$ws.Cells.Item(6, 11).Value = 62
$ws.Cells.Item(7, 11).Value = 179

$ws = $wb.Worksheets.Item(65)
$ws.Cells.Item(4, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 161

$ws = $wb.Worksheets.Item(68)
$ws.Cells.Item(2, 11).Value = 55
$ws.Cells.Item(6, 11).Value = 91
$ws.Cells.Item(7, 11).Value = 219

$ws = $wb.Worksheets.Item(71)
$ws.Cells.Item(6, 11).Value = 69
$ws.Cells.Item(7, 11).Value = 189

$ws = $wb.Worksheets.Item(74)
$ws.Cells.Item(6, 11).Value = 44
$ws.Cells.Item(7, 11).Value = 187

$ws = $wb.Worksheets.Item(75)
$ws.Cells.Item(2, 11).Value = 74
$ws.Cells.Item(3, 11).Value = 68
$ws.Cells.Item(6, 11).Value = 83
$ws.Cells.Item(7, 11).Value = 259

$ws = $wb.Worksheets.Item(76)
$ws.Cells.Item(6, 11).Value = 14
$ws.Cells.Item(7, 11).Value = 54

$ws = $wb.Worksheets.Item(82)
$ws.Cells.Item(3, 11).Value = 25
$ws.Cells.Item(6, 11).Value = 51
$ws.Cells.Item(7, 11).Value = 97

$ws = $wb.Worksheets.Item(84)
$ws.Cells.Item(2, 11).Value = 62
$ws.Cells.Item(7, 11).Value = 143

$ws = $wb.Worksheets.Item(85)
$ws.Cells.Item(3, 11).Value = 9
$ws.Cells.Item(7, 11).Value = 26

$ws = $wb.Worksheets.Item(87)
$ws.Cells.Item(2, 11).Value = 18
$ws.Cells.Item(7, 11).Value = 71

$ws = $wb.Worksheets.Item(91)
$ws.Cells.Item(2, 11).Value = 10
$ws.Cells.Item(7, 11).Value = 36

$ws = $wb.Worksheets.Item(96)
$ws.Cells.Item(2, 11).Value = 6
$ws.Cells.Item(7, 11).Value = 15
